# Generate Report for Handoff
# Regenerates the localization-status report: new source file UUID
# (f66dac40-... -> 9b43d845-...), new xliff content hash
# (34e2cda2... -> e7c78d2b...), refreshed handoff timestamps, and the
# "Latest Target File" / "Latest Handback File" columns reset because the
# handback step hasn't completed yet for this run (Latest Handback
# DateTime reset to the zero-date sentinel 0001-01-01 00:00:00).

$wb = $excel.ActiveWorkbook

$oldGuid = "f66dac40-2ba3-46b8-a18d-b957bef67687"
$newGuid = "9b43d845-f86e-4d51-a78c-c8f3c04e59dd"

$newFileName      = "$newGuid.md"
$newPathAndName   = "e2e\$newGuid.md"
$newZhCnXlf       = "$newGuid.e7c78d2b99f8de0fd1c31e39b5f68e621b3ca522.zh-cn.xlf"
$newDeDeXlf       = "$newGuid.e7c78d2b99f8de0fd1c31e39b5f68e621b3ca522.de-de.xlf"

$newHoXliffDate   = "2016-08-17 22:58:19"
$newZhCnHandoffDt = "2016-08-17 22:58:14"
$zeroDateTime     = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = $newPathAndName
    }
}
$wsOverview.Range("G2").Value = $newHoXliffDate

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFileName
foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFileName
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHandoffDt
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = $zeroDateTime

$wsZhCn.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsZhCn.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFileName
foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFileName
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newHoXliffDate
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = $zeroDateTime

$wsDeDe.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsDeDe.Columns.Item(10).ColumnWidth = 20.833333333333332
